# Updates cryptos list values (Price column D, Volume(1h) column E)
# Mirrors the commit "Updated cryptos list on Sat Aug  3 09:57:59 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.565.08"
$ws.Range("E2").Value = "  -4.27%  "
$ws.Range("D3").Value = "2.975.73"
$ws.Range("E3").Value = "  -5.18%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.00%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -1.11%  "
$ws.Range("D9").Value = "2.984.93"
$ws.Range("E9").Value = "  -5.29%  "
$ws.Range("E10").Value = "  -3.67%  "
$ws.Range("E11").Value = "  -7.10%  "
$ws.Range("E12").Value = "  -4.83%  "
$ws.Range("D13").Value = "3.494.50"
$ws.Range("E13").Value = "  -5.27%  "
$ws.Range("E14").Value = "  -1.71%  "
$ws.Range("D15").Value = "61.632.76"
$ws.Range("E15").Value = "  -4.24%  "
$ws.Range("E16").Value = "  -5.64%  "
$ws.Range("D17").Value = "2.978.56"
$ws.Range("E17").Value = "  -5.17%  "
$ws.Range("E18").Value = "  -5.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.01%  "
$ws.Range("E20").Value = "  -4.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.55%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  -3.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.41%  "
$ws.Range("E26").Value = "  -2.98%  "
$ws.Range("D27").Value = "3.097.02"
$ws.Range("E27").Value = "  -5.45%  "
$ws.Range("E28").Value = "  -3.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").Value = "0.0₃0943"
$ws.Range("E30").Value = "  -7.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.83%  "
$ws.Range("E33").Value = "  -4.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.46"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "159.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.08%  "
$ws.Range("E36").Value = "  -4.78%  "
$ws.Range("E37").Value = "  -5.97%  "
$ws.Range("E38").Value = "  -3.64%  "
$ws.Range("E39").Value = "  -5.77%  "
$ws.Range("E40").Value = "  -8.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E43").Value = "  -8.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.90%  "
$ws.Range("E45").Value = "  -2.66%  "
$ws.Range("E46").Value = "  -3.46%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("E48").Value = "  -6.39%  "
$ws.Range("E49").Value = "  -3.57%  "
$ws.Range("E50").Value = "  -2.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.69%  "
